# Update New Zealand MSME summary figures to higher-precision values.
# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (these are stored as text, not numbers, in
# the original workbook, so we must force a text number format before
# assigning the value and then restore the cell's original style so the
# on-disk style index is unaffected).
$updates = @{
    "B11" = "21.97"  # Enterprises density (per 1000 people) - Micro
    "C11" = "10.68"  # Enterprises density (per 1000 people) - SMEs
    "D11" = "32.65"  # Enterprises density (per 1000 people) - MSMEs
    "B12" = "11.64"  # Employment (% of total) - Micro
    "C12" = "41.19"  # Employment (% of total) - SMEs
    "D12" = "52.83"  # Employment (% of total) - MSMEs
    "B14" = "66.32"  # Enterprises (% of total) - Micro
    "C14" = "32.23"  # Enterprises (% of total) - SMEs
    "D14" = "98.56"  # Enterprises (% of total) - MSMEs
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = $originalStyle
}
